$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3657440749772013
$ws.Range("D2").Value = 0.7180493144005826

$ws.Range("C3").Value = -0.2290136610879999
$ws.Range("D3").Value = 0.8209738576604484

$ws.Range("C4").Value = -0.05391290380741177
$ws.Range("D4").Value = 0.9574911437805844

$ws.Range("C5").Value = -0.637776099894246
$ws.Range("D5").Value = 0.5302021728374311

$ws.Range("C6").Value = -0.4787898061666366
$ws.Range("D6").Value = 0.6368124824370343

$ws.Range("C7").Value = -0.2783401071444002
$ws.Range("D7").Value = 0.7833506196637554

$ws.Range("C8").Value = -0.7165546760659749
$ws.Range("D8").Value = 0.4811871872133406

$ws.Range("C9").Value = 0.1306181287354911
$ws.Range("D9").Value = 0.8972644317933072

$ws.Range("C10").Value = -0.2927998974778782
$ws.Range("D10").Value = 0.772418964966223

$ws.Range("C11").Value = -0.4224565701322059
$ws.Range("D11").Value = 0.6767920765556044
